# Applies the commit:
#   "Diseno de pruebas para busqueda, ordenamiento y adicion de jugadores
#    Ahora los jugadores se ordenan por Nombre en el arbol binario"
#
# Target sheet: "Must Have" (sheet2.xml)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Must Have")

# --- Row 3: CompareTo cell (E3) gets the highlight fill style (same text) ---
$ws.Range("E3").Interior.ThemeColor = 9
$ws.Range("E3").Interior.TintAndShade = 0.59999389629810485

# --- Row 4: "Por nombre" cell (C4) now highlighted with the NEW accent color ---
$ws.Range("C4").Interior.ThemeColor = 5
$ws.Range("C4").Interior.TintAndShade = 0.39997558519241921

# --- Row 4: add a note cell H4 ---
$ws.Range("H4").Value = "vale si es de la manera recursiva?"

# --- Row 5: add new test-case name in C5, with the same highlight used elsewhere ---
$ws.Range("C5").Value = "buscarJugadorPuntos"
$ws.Range("C5").Interior.ThemeColor = 9
$ws.Range("C5").Interior.TintAndShade = 0.59999389629810485

# --- Row 9: rename exception test cases, highlight C9/D9/E9 ---
$ws.Range("C9").Value = "NombreNoExiste"
$ws.Range("C9").Interior.ThemeColor = 9
$ws.Range("C9").Interior.TintAndShade = 0.59999389629810485

$ws.Range("D9").Interior.ThemeColor = 9
$ws.Range("D9").Interior.TintAndShade = 0.59999389629810485

$ws.Range("E9").Value = "JugadorRepetidoException"
$ws.Range("E9").Interior.ThemeColor = 9
$ws.Range("E9").Interior.TintAndShade = 0.59999389629810485

# --- Update selections to match the saved view state ---
$ws1 = $wb.Worksheets.Item("Requerimientos")
$ws1.Range("A8").Select()

$ws.Range("D10").Select()
